# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the FFXIV leve-profit workbook.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 108.625
$ws.Range("I8").Value = 108.625
$ws.Range("K8").Value = 325.875
$ws.Range("M8").Value = -186.875
$ws.Range("H9").Value = 539.2353000000001
$ws.Range("I9").Value = 720.0909
$ws.Range("J9").Value = 207.66667
$ws.Range("K9").Value = 720.0909
$ws.Range("L9").Value = 207.66667
$ws.Range("M9").Value = -551.0909
$ws.Range("N9").Value = -545.6666700000001
$ws.Range("H69").Value = 19974.75
$ws.Range("I69").Value = 19898
$ws.Range("K69").Value = 59694
$ws.Range("M69").Value = -58820
$ws.Range("H72").Value = 19974.75
$ws.Range("I72").Value = 19898
$ws.Range("K72").Value = 179082
$ws.Range("M72").Value = -174714
$ws.Range("H113").Value = 6111.364
$ws.Range("I113").Value = 4404.6665
$ws.Range("J113").Value = 8159.4
$ws.Range("K113").Value = 4404.6665
$ws.Range("L113").Value = 8159.4
$ws.Range("M113").Value = -1150.6665
$ws.Range("N113").Value = -14667.4
$ws.Range("H116").Value = 7876.125
$ws.Range("I116").Value = 4919.6
$ws.Range("K116").Value = 4919.6
$ws.Range("M116").Value = -1477.6
$ws.Range("H132").Value = 6079.3335
$ws.Range("I132").Value = 6193.1577
$ws.Range("K132").Value = 18579.4731
$ws.Range("M132").Value = -16049.4731

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1822.2
$ws.Range("I32").Value = 1773.7858
$ws.Range("K32").Value = 1773.7858
$ws.Range("M32").Value = -1486.7858
$ws.Range("H68").Value = 59999
$ws.Range("J68").Value = 59999
$ws.Range("L68").Value = 59999
$ws.Range("N68").Value = -61621
$ws.Range("H71").Value = 59999
$ws.Range("J71").Value = 59999
$ws.Range("L71").Value = 179997
$ws.Range("N71").Value = -188109
$ws.Range("H74").Value = 1761.5
$ws.Range("I74").Value = 1471.9166
$ws.Range("K74").Value = 1471.9166
$ws.Range("M74").Value = -597.9166
$ws.Range("H77").Value = 1761.5
$ws.Range("I77").Value = 1471.9166
$ws.Range("K77").Value = 7359.583000000001
$ws.Range("M77").Value = -2991.583000000001
$ws.Range("H97").Value = 499.75
$ws.Range("I97").Value = 564
$ws.Range("J97").Value = 307
$ws.Range("K97").Value = 564
$ws.Range("L97").Value = 307
$ws.Range("M97").Value = -68
$ws.Range("N97").Value = -1299

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 945
$ws.Range("I20").Value = 1043.9231
$ws.Range("J20").Value = 516.3333
$ws.Range("K20").Value = 1043.9231
$ws.Range("L20").Value = 516.3333
$ws.Range("M20").Value = -796.9231
$ws.Range("N20").Value = -1010.3333
$ws.Range("H94").Value = 394.2
$ws.Range("I94").Value = 371.33334
$ws.Range("J94").Value = 600
$ws.Range("K94").Value = 371.33334
$ws.Range("L94").Value = 600
$ws.Range("M94").Value = 79.66665999999998
$ws.Range("N94").Value = -1502

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4795.5
$ws.Range("I31").Value = 4795.5
$ws.Range("K31").Value = 4795.5
$ws.Range("M31").Value = -4500.5
$ws.Range("H34").Value = 4795.5
$ws.Range("I34").Value = 4795.5
$ws.Range("K34").Value = 4795.5
$ws.Range("M34").Value = -4593.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 87.125
$ws.Range("I2").Value = 63.5
$ws.Range("J2").Value = 110.75
$ws.Range("K2").Value = 381
$ws.Range("L2").Value = 664.5
$ws.Range("M2").Value = -268
$ws.Range("N2").Value = -890.5
$ws.Range("H12").Value = 158.33333
$ws.Range("I12").Value = 188.83333
$ws.Range("K12").Value = 566.49999
$ws.Range("M12").Value = -393.49999
$ws.Range("H68").Value = 3150
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3150
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 9450
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -11072
$ws.Range("H71").Value = 3150
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3150
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 28350
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -36462

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 235.38461
$ws.Range("I2").Value = 239.75
$ws.Range("J2").Value = 228.4
$ws.Range("K2").Value = 239.75
$ws.Range("L2").Value = 228.4
$ws.Range("M2").Value = -126.75
$ws.Range("N2").Value = -454.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1858.8334
$ws.Range("I22").Value = 1599.1428
$ws.Range("J22").Value = 2024.091
$ws.Range("K22").Value = 1599.1428
$ws.Range("L22").Value = 2024.091
$ws.Range("M22").Value = -1304.1428
$ws.Range("N22").Value = -2614.091
$ws.Range("H27").Value = 1858.8334
$ws.Range("I27").Value = 1599.1428
$ws.Range("J27").Value = 2024.091
$ws.Range("K27").Value = 1599.1428
$ws.Range("L27").Value = 2024.091
$ws.Range("M27").Value = -1492.1428
$ws.Range("N27").Value = -2238.091
$ws.Range("H33").Value = 3000
$ws.Range("I33").Value = 3000
$ws.Range("K33").Value = 3000
$ws.Range("M33").Value = -2710
$ws.Range("H61").Value = 2219.4
$ws.Range("I61").Value = 1788.75
$ws.Range("K61").Value = 1788.75
$ws.Range("M61").Value = -1586.75
$ws.Range("H100").Value = 2929.5
$ws.Range("I100").Value = 2929.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2929.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2388.5
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 2219.4
$ws.Range("I113").Value = 1788.75
$ws.Range("K113").Value = 1788.75
$ws.Range("M113").Value = 381.25
$ws.Range("H132").Value = 2161.318
$ws.Range("I132").Value = 2107.55
$ws.Range("J132").Value = 2699
$ws.Range("K132").Value = 6322.650000000001
$ws.Range("L132").Value = 8097
$ws.Range("M132").Value = -3792.650000000001
$ws.Range("N132").Value = -13157

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6333.8335
$ws.Range("I122").Value = 6800.6
$ws.Range("K122").Value = 20401.8
$ws.Range("M122").Value = -17951.8
$ws.Range("H132").Value = 3851.5715
$ws.Range("I132").Value = 1326
$ws.Range("J132").Value = 19005
$ws.Range("K132").Value = 3978
$ws.Range("L132").Value = 57015
$ws.Range("M132").Value = -1448
$ws.Range("N132").Value = -62075
$ws.Range("H136").Value = 2971.1052
$ws.Range("I136").Value = 2597
$ws.Range("K136").Value = 7791
$ws.Range("M136").Value = -5241
